$d = $word.ActiveDocument

$pairs = @(
    @{ old = "524×8="; new = "941×7=" },
    @{ old = "572×3="; new = "254×9=" },
    @{ old = "263×5="; new = "296×4=" },
    @{ old = "846×6="; new = "999×7=" },
    @{ old = "363×8="; new = "827×8=" },
    @{ old = "734×2="; new = "102×6=" },
    @{ old = "992×9="; new = "242×2=" },
    @{ old = "989×9="; new = "866×2=" },
    @{ old = "747×4="; new = "428×2=" },
    @{ old = "432×2="; new = "883×7=" },
    @{ old = "895×7="; new = "955×8=" },
    @{ old = "985×6="; new = "482×2=" },
    @{ old = "389×7="; new = "468×8=" },
    @{ old = "428×9="; new = "315×6=" },
    @{ old = "386×5="; new = "165×2=" },
    @{ old = "751×3="; new = "769×5=" },
    @{ old = "971×9="; new = "622×4=" },
    @{ old = "657×7="; new = "524×3=" },
    @{ old = "856×9="; new = "482×5=" },
    @{ old = "389×3="; new = "889×3=" },
    @{ old = "444×9="; new = "634×2=" },
    @{ old = "261×2="; new = "825×8=" },
    @{ old = "475×2="; new = "525×9=" },
    @{ old = "873×8="; new = "666×9=" },
    @{ old = "146×5="; new = "370×9=" }
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
